$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells keep their exact string representation (matching the
# original inline-string/text cell type) instead of being auto-converted to numbers.
function Set-TextCell($sheet, $addr, $val) {
    $c = $sheet.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
}

Set-TextCell $ws "D2" "243.11"
Set-TextCell $ws "D3" "23.01"
$ws.Range("B4").Value = "HuobiToken"
$ws.Range("C4").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextCell $ws "D4" "5.394"
$ws.Range("E4").Value = "3HuobiTokenHT"
$ws.Range("B5").Value = "Cronos"
$ws.Range("C5").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextCell $ws "D5" "0.05934"
$ws.Range("E5").Value = "4CronosCRO"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell $ws "D6" "3.393"
$ws.Range("E6").Value = "5GateTokenGT"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextCell $ws "D7" "6.464"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws "D8" "0.8080"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell $ws "D9" "0.9053"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell $ws "D10" "0.1417"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell $ws "D11" "0.07444"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell $ws "D12" "0.03279"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell $ws "D13" "0.03044"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell $ws "D14" "0.09325"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell $ws "D15" "3.943"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell $ws "D16" "0.001574"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell $ws "D17" "0.04788"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell $ws "D18" "0.0005944"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("B19").Value = "TigerCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell $ws "D19" "0.006197"
$ws.Range("E19").Value = "18TigerCashTCH"
$ws.Range("B20").Value = "UpBots"
$ws.Range("C20").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextCell $ws "D20" "0.007493"
$ws.Range("E20").Value = "19UpBotsUBXTBestin24h"
$ws.Range("B21").Value = "HotbitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell $ws "D21" "0.004413"
$ws.Range("E21").Value = "20HotbitTokenHTB"
$ws.Range("B22").Value = "BitKan"
$ws.Range("C22").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell $ws "D22" "0.0009879"
$ws.Range("E22").Value = "21BitKanKAN"
$ws.Range("B23").Value = "NitroEx"
$ws.Range("C23").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextCell $ws "D23" "0.00007806"
$ws.Range("E23").Value = "22NitroExNTX"
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D24" "3.615"
$ws.Range("E24").Value = "23LEOLEO"
Set-TextCell $ws "D40" "0.03869"
Set-TextCell $ws "D41" "0.006218"
Set-TextCell $ws "D42" "0.1068"
Set-TextCell $ws "D43" "0.002612"
Set-TextCell $ws "D44" "0.007255"
Set-TextCell $ws "D45" "0.00005197"
Set-TextCell $ws "D46" "0.00000000751"
Set-TextCell $ws "D47" "0.0005804"
Set-TextCell $ws "D48" "0.9607"
Set-TextCell $ws "D50" "0.00002102"
Set-TextCell $ws "D51" "0.0002002"
